$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "94.435.35"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.076.51"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.62"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "610.56"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.10"
$ws.Range("E7").Value = "  +1.77%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.377"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.808"
$ws.Range("E10").Value = "  +10.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.073.19"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "94.073.32"
$ws.Range("E13").Value = "  +1.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000241"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.82"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.35"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.641.84"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.061.89"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.56"
$ws.Range("E19").Value = "  -4.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.38"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.67"
$ws.Range("E21").Value = "  -0.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "439.78"
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.81"
$ws.Range("E23").Value = "  -5.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000189"
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.33"
$ws.Range("E25").Value = "  +6.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.52"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "84.71"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.83"
$ws.Range("E28").Value = "  +2.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.236.34"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.249"
$ws.Range("E31").Value = "  +10.22%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.176"
$ws.Range("E32").Value = "  +5.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.124"
$ws.Range("E33").Value = "  -4.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.01"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.988"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.68"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("E37").Value = "  -1.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.40"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.06"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "477.02"
$ws.Range("E41").Value = "  -0.68%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.437"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.74"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("E46").Value = "  -5.24%  "
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.672"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("E49").Value = "  -2.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.59"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").Value = "  +0.25%  "